$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tabella 13")

# Update row 15: the old "di Stato di famiglia per iscritto in convivenza anagrafica"
# label is replaced with "di Residenza in convivenza".
$ws.Range("B15").Value = "di Residenza in convivenza"

# Add a new row 21 with id 20 and the new certificate type label.
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Storico di Residenza alla data"

# Update the active selection to match the authored workbook state.
$ws.Range("B21").Select()
